$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

# --- Expand the table to include a new 7th column (appended at the end as "G") ---
$tbl.Resize($ws.Range("A1:G1048576"))

# --- Shift the old columns C..F (Dish, Notes, Assigned to, Completed) right by one,
#     into D..G, freeing up column C for the new "Elapsed Time (min)" column. ---
$dishHeader      = $ws.Cells.Item(1,3).Value2
$notesHeader     = $ws.Cells.Item(1,4).Value2
$assignedHeader  = $ws.Cells.Item(1,5).Value2
$completedHeader = $ws.Cells.Item(1,6).Value2

$ws.Cells.Item(1,7).Value = $completedHeader
$ws.Cells.Item(1,6).Value = $assignedHeader
$ws.Cells.Item(1,5).Value = $notesHeader
$ws.Cells.Item(1,4).Value = $dishHeader
$ws.Cells.Item(1,3).Value = "Elapsed Time (min)"

for ($r = 2; $r -le 5; $r++) {
    $dishV  = $ws.Cells.Item($r,3).Value2
    $notesV = $ws.Cells.Item($r,4).Value2

    $ws.Cells.Item($r,5).Value = $notesV
    $ws.Cells.Item($r,4).Value = $dishV
    $ws.Cells.Item($r,3).ClearContents()
}

# --- Populate the new "Elapsed Time (min)" column (C) with its data ---
$ws.Cells.Item(2,3).Value = 25
$ws.Cells.Item(3,3).Value = 5
$ws.Cells.Item(4,3).Value = 10
$ws.Cells.Item(5,3).Value = "-"

# --- Populate "Assigned to" column (F) ---
$ws.Cells.Item(2,6).Value = "Frank"
$ws.Cells.Item(3,6).Value = "Sam"
$ws.Cells.Item(4,6).Value = "Charlie"
$ws.Cells.Item(5,6).Value = "-"

# --- Center-align the "Prep Time (Mins)" (B) and "Elapsed Time (min)" (C) columns,
#     matching the dxf/cellXfs the table keeps for those two data columns. ---
$ws.Range("B1:B5").HorizontalAlignment = -4108
$ws.Range("C1:C5").HorizontalAlignment = -4108

# --- Column widths / bestfit to roughly match the edited layout ---
$ws.Columns.Item(2).ColumnWidth = 17.75
$ws.Columns.Item(3).ColumnWidth = 17.75
$ws.Columns.Item(2).AutoFit()

# --- Selection, matching the author's last cursor position ---
$ws.Range("F15").Select()
